$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "maa://24702 (94.1), maa://25390 (97.42), maa://36681 (92.06)"
$ws.Range("K3").Value = "*maa://22880 (70.2), maa://20276 (82.48), *maa://22749 (62.5)"
$ws.Range("O3").Value = "maa://21249 (95.57), maa://26254 (95.24)"
$ws.Range("W4").Value = "**maa://32495 (48.12), ***maa://31785 (15.74), ***maa://36683 (26.67)"
$ws.Range("K6").Value = "maa://24839 (99.2)"
$ws.Range("G7").Value = "*maa://22763 (65.38)"
$ws.Range("AE7").Value = "*maa://26191 (69.44), *maa://36671 (73.81)"
$ws.Range("W8").Value = "maa://21411 (96.28)"
$ws.Range("W9").Value = "maa://26223 (96.88)"
$ws.Range("C10").Value = "***maa://25695 (19.41), **maa://32237 (38.89), ***maa://34206 (14.29), ***maa://39951 (20.0), ***maa://39243 (25.0)"
$ws.Range("S10").Value = "maa://27395 (97.26), maa://22755 (87.5), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("W12").Value = "maa://22753 (91.72), *maa://21485 (76.56), maa://37962 (81.25)"
$ws.Range("AA12").Value = "maa://23669 (95.83), maa://36677 (94.74), maa://39872 (81.82)"
$ws.Range("O13").Value = "maa://22676 (91.75), *maa://22583 (74.58), *maa://22500 (55.81)"
$ws.Range("W13").Value = "*maa://34957 (77.5), *maa://22768 (53.33)"
$ws.Range("G15").Value = "maa://24304 (88.33), maa://21478 (91.18)"
$ws.Range("C16").Value = "maa://21441 (96.15), maa://36679 (93.33), maa://37650 (95.24)"
$ws.Range("AA16").Value = "maa://26228 (96.05)"
$ws.Range("G17").Value = "maa://22430 (88.57), *maa://39599 (80.0)"
$ws.Range("G18").Value = "maa://24421 (90.38)"
$ws.Range("S19").Value = "maa://24386 (98.72)"
$ws.Range("K20").Value = "maa://41331 (86.36)"
$ws.Range("W21").Value = "maa://20110 (86.57), maa://34946 (90.32)"
$ws.Range("AA21").Value = "*maa://21443 (78.7), ***maa://23820 (29.63)"
$ws.Range("AE21").Value = "maa://22524 (94.25), *maa://22432 (75.47)"
$ws.Range("B23").Value = "'2"
$ws.Range("C23").Value = "***maa://28036 (28.79), maa://41753 (100.0)"
$ws.Range("K23").Value = "maa://39756 (92.14), maa://39875 (95.0)"
$ws.Range("W23").Value = "*maa://28503 (60.71)"
$ws.Range("C24").Value = "maa://24368 (80.43)"
$ws.Range("W24").Value = "maa://23504 (93.04), maa://29988 (85.86), **maa://22892 (40.43), *maa://25141 (76.86), *maa://36663 (79.25), ***maa://22815 (23.08)"
$ws.Range("AE24").Value = "maa://22523 (84.86), *maa://36672 (75.61), maa://29910 (93.88), **maa://21440 (34.55)"
$ws.Range("G25").Value = "*maa://29063 (77.34), *maa://25311 (74.73), ***maa://22725 (4.84)"
$ws.Range("G27").Value = "**maa://21283 (49.32), maa://34494 (100.0), **maa://36665 (44.44), maa://39601 (85.71)"
$ws.Range("S27").Value = "*maa://30624 (77.78)"
$ws.Range("V28").Value = "'3"
$ws.Range("W28").Value = "maa://39929 (85.71), ***maa://39723 (15.15), maa://41749 (100.0)"
$ws.Range("C29").Value = "maa://31694 (97.78)"
$ws.Range("AE29").Value = "*maa://24080 (68.33), ***maa://34960 (9.09)"
$ws.Range("S32").Value = "maa://41108 (91.67), maa://41238 (93.94)"
$ws.Range("K35").Value = "maa://41296 (96.43)"
$ws.Range("AE38").Value = "maa://36697 (86.18)"
$ws.Range("O40").Value = "maa://23278 (96.19), maa://21386 (95.63), maa://36664 (92.31)"
$ws.Range("O41").Value = "**maa://35616 (36.0)"
$ws.Range("G44").Value = "maa://29768 (97.49), maa://27728 (96.0)"
